$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are plain text (names/URLs); D holds numeric-looking
# strings that must remain TEXT (preserve formats like "1.00", "0.580",
# multi-dot big numbers, leading zeros, etc.); E holds percentage strings
# that are already safe as text because of the surrounding spaces and %.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '96.275.60'
$ws.Range('E2').Value = '  -0.79%  '
Set-TextValue $ws.Range('D3') '3.460.07'
$ws.Range('E3').Value = '  +4.00%  '
$ws.Range('E4').Value = '  +0.39%  '
Set-TextValue $ws.Range('D5') '244.95'
$ws.Range('E5').Value = '  -0.82%  '
Set-TextValue $ws.Range('D6') '645.84'
$ws.Range('E6').Value = '  -0.81%  '
Set-TextValue $ws.Range('D7') '1.42'
$ws.Range('E7').Value = '  +4.37%  '
Set-TextValue $ws.Range('D8') '0.409'
$ws.Range('E8').Value = '  -0.80%  '
Set-TextValue $ws.Range('D9') '1.00'
$ws.Range('E9').Value = '  +0.22%  '
Set-TextValue $ws.Range('D10') '0.998'
$ws.Range('E10').Value = '  +1.94%  '
Set-TextValue $ws.Range('D11') '3.456.36'
$ws.Range('E11').Value = '  +3.95%  '
Set-TextValue $ws.Range('D12') '43.70'
$ws.Range('E12').Value = '  +9.84%  '
Set-TextValue $ws.Range('D13') '0.201'
$ws.Range('E13').Value = '  -1.55%  '
Set-TextValue $ws.Range('D14') '6.14'
$ws.Range('E14').Value = '  +2.86%  '
Set-TextValue $ws.Range('D15') '96.448.11'
$ws.Range('E15').Value = '  -0.40%  '
Set-TextValue $ws.Range('D16') '4.118.60'
$ws.Range('E16').Value = '  +4.52%  '
Set-TextValue $ws.Range('D17') '0.0000252'
$ws.Range('E17').Value = '  +0.96%  '
Set-TextValue $ws.Range('D18') '8.65'
$ws.Range('E18').Value = '  +1.47%  '
Set-TextValue $ws.Range('D19') '3.433.14'
$ws.Range('E19').Value = '  +3.41%  '
Set-TextValue $ws.Range('D20') '18.25'
$ws.Range('E20').Value = '  +9.38%  '
Set-TextValue $ws.Range('D21') '11.78'
$ws.Range('E21').Value = '  +12.95%  '
Set-TextValue $ws.Range('D22') '0.500'
$ws.Range('E22').Value = '  +4.88%  '
Set-TextValue $ws.Range('D23') '516.04'
$ws.Range('E23').Value = '  +5.30%  '
Set-TextValue $ws.Range('D24') '3.29'
$ws.Range('E24').Value = '  +0.98%  '
Set-TextValue $ws.Range('D25') '0.0000195'
$ws.Range('E25').Value = '  -0.61%  '
Set-TextValue $ws.Range('D26') '6.70'
$ws.Range('E26').Value = '  +5.53%  '
Set-TextValue $ws.Range('D27') '92.58'
$ws.Range('E27').Value = '  +0.40%  '
Set-TextValue $ws.Range('D28') '12.38'
$ws.Range('E28').Value = '  +3.80%  '
Set-TextValue $ws.Range('D29') '3.669.47'
$ws.Range('E29').Value = '  +4.96%  '
Set-TextValue $ws.Range('D30') '11.98'
$ws.Range('E30').Value = '  +11.61%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D31') '1.00'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D32') '2.76'
$ws.Range('E32').Value = '  +12.40%  '
Set-TextValue $ws.Range('D33') '0.140'
$ws.Range('E33').Value = '  -1.20%  '
Set-TextValue $ws.Range('D34') '0.185'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range('D35') '0.580'
$ws.Range('E35').Value = '  +7.08%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D36') '30.77'
$ws.Range('E36').Value = '  +10.62%  '
Set-TextValue $ws.Range('D37') '0.999'
$ws.Range('E37').Value = '  -0.33%  '
Set-TextValue $ws.Range('D38') '7.83'
$ws.Range('E38').Value = '  +4.79%  '
Set-TextValue $ws.Range('D39') '1.46'
$ws.Range('E39').Value = '  +0.12%  '
Set-TextValue $ws.Range('D40') '0.152'
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('E41').Value = '  +0.02%  '
Set-TextValue $ws.Range('D42') '506.84'
$ws.Range('E42').Value = '  +1.19%  '
Set-TextValue $ws.Range('D43') '0.903'
$ws.Range('E43').Value = '  +10.09%  '
Set-TextValue $ws.Range('D44') '24.25'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D45') '1.71'
$ws.Range('E45').Value = '  +5.40%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D46') '0.0419'
$ws.Range('E46').Value = '  +3.75%  '
Set-TextValue $ws.Range('D47') '3.35'
$ws.Range('E47').Value = '  +7.82%  '
Set-TextValue $ws.Range('D48') '5.55'
$ws.Range('E48').Value = '  +2.94%  '
$ws.Range('B49').Value = 'MantraDAO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws.Range('D49') '3.58'
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D50') '2.20'
$ws.Range('E50').Value = '  +12.32%  '
Set-TextValue $ws.Range('D51') '8.41'
$ws.Range('E51').Value = '  +1.18%  '
